# Update the "Förändrad" (Changed) date column (C) for rows 2-45:
# increment the stored date serial value by 1 day (2023-10-03 -> 2023-10-04)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
